$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title
#    paragraph: an empty leading run, a bold "Meta description" run, and
#    a plain run with the rest of the description text.
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphAfter()

$metaP = $d.Paragraphs(2)
$metaP.Style = "Normal"

$metaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Explore the innovative features of Avatar: Gateway Guardians slot game and play for free. Read our review for a unique playing experience.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$metaP.Range.InsertXML($metaXml)

# ---------------------------------------------------------------------
# 2) Remove the duplicated bold "Play Avatar: Gateway Guardians Slot for
#    Free - Review" paragraph that used to sit near the end of the doc
#    (right before the italic image-prompt paragraph).
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$dupTitleP = $d.Paragraphs($count - 1)
$dupTitleP.Range.Delete()

# ---------------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph with the new
#    image-generation prompt, preserving the italic run formatting and
#    the leading empty run.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastP = $d.Paragraphs($count)

# Clear the paragraph's content but keep its trailing paragraph mark.
$contentRng = $d.Range($lastP.Range.Start, $lastP.Range.End - 1)
$contentRng.Delete()

$count = $d.Paragraphs.Count
$lastP = $d.Paragraphs($count)

$newPrompt = 'Please create a feature image for "Avatar: Gateway Guardians" that fits the following criteria: - It should be in a cartoon style. - The main character in the image should be a happy Maya warrior wearing glasses. The image should feature a round frame, similar to the circular reels in the game. The happy Maya warrior should be standing in the center of the frame wielding a staff adorned with blue and azure hues. The warrior should be wearing a traditional Mayan headdress and gray armor. The background should be misty, with a few floating rocks and a hint of blue and green hues. The overall tone of the image should be vibrant and exciting, inviting players to try their luck in the world of Avatar: Gateway Guardians.'

$promptXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>' + $newPrompt + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$countBeforeInsert = $d.Paragraphs.Count
$lastP.Range.InsertXML($promptXml)
$countAfterInsert = $d.Paragraphs.Count

# InsertXML on the very last paragraph of the body leaves behind a stray
# empty paragraph (the original trailing paragraph mark survives as its
# own empty paragraph). Merge it back by deleting the paragraph mark that
# now separates our new content paragraph from that stray empty one.
if ($countAfterInsert -gt $countBeforeInsert) {
    $contentP = $d.Paragraphs($countAfterInsert - 1)
    $markRng = $d.Range($contentP.Range.End - 1, $contentP.Range.End)
    $markRng.Delete()
}

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
